$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Fitness column (C) values.
# Row 2 (Generation 0) drops to 7295, all subsequent rows (Generation 1-250)
# drop to 7293.
$ws.Range("C2").Value = 7295
$ws.Range("C3:C252").Value = 7293
